# The document ends with a run of empty "Caption"-styled paragraphs
# (used as spacing/placeholder paragraphs after the table of contents).
# The edit trims that run down to just the first two, deleting the
# trailing six empty Caption paragraphs that sit right before sectPr.

$d = $word.ActiveDocument

# Locate every empty paragraph styled "Caption", in document order.
$captionIndices = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Caption" -and $p.Range.Text.Length -le 1) {
        $captionIndices += $i
    }
}

# Keep the first two empty Caption paragraphs; delete the remaining
# (trailing) ones as a single contiguous range.
$keep = 2
if ($captionIndices.Count -gt $keep) {
    $firstToDelete = $captionIndices[$keep]
    $lastToDelete = $captionIndices[$captionIndices.Count - 1]

    $rangeStart = $d.Paragraphs.Item($firstToDelete).Range.Start
    $rangeEnd = $d.Paragraphs.Item($lastToDelete).Range.End

    $d.Range($rangeStart, $rangeEnd).Delete()
}
